# Apply Phantom_Profits value corrections (scheduled-runner refresh of pricing data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 22552.334
$ws.Range("J3").Value = 22552.334
$ws.Range("L3").Value = 22552.334
$ws.Range("N3").Value = -22780.334

$ws.Range("H33").Value = 338
$ws.Range("I33").Value = 350.72726
$ws.Range("K33").Value = 350.72726
$ws.Range("M33").Value = -121.72726

$ws.Range("H43").Value = 7513.75
$ws.Range("J43").Value = 7351.3335
$ws.Range("L43").Value = 7351.3335
$ws.Range("N43").Value = -7489.3335

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

$ws.Range("H96").Value = 1283.0769
$ws.Range("J96").Value = 433.8
$ws.Range("L96").Value = 1301.4
$ws.Range("N96").Value = -4047.4

$ws.Range("H102").Value = 22552.334
$ws.Range("J102").Value = 22552.334
$ws.Range("L102").Value = 22552.334
$ws.Range("N102").Value = -29042.334

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws.Range("H131").Value = 990.5
$ws.Range("I131").Value = 990.5
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2971.5
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2068.5
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 4168.32
$ws.Range("I132").Value = 4055.318
$ws.Range("K132").Value = 12165.954
$ws.Range("M132").Value = -9635.954000000002

$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

$ws.Range("H135").Value = 817.7
$ws.Range("I135").Value = 353.33334
$ws.Range("J135").Value = 4997
$ws.Range("K135").Value = 3180.00006
$ws.Range("L135").Value = 44973
$ws.Range("M135").Value = -645.0000600000003
$ws.Range("N135").Value = -50043

$ws.Range("H136").Value = 75500
$ws.Range("J136").Value = 75500
$ws.Range("L136").Value = 75500
$ws.Range("N136").Value = -85700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5027.1177
$ws.Range("I61").Value = 5162.2666
$ws.Range("J61").Value = 4013.5
$ws.Range("K61").Value = 5162.2666
$ws.Range("L61").Value = 4013.5
$ws.Range("M61").Value = -4950.2666
$ws.Range("N61").Value = -4437.5

$ws.Range("H94").Value = 27876.334
$ws.Range("J94").Value = 27876.334
$ws.Range("L94").Value = 27876.334
$ws.Range("N94").Value = -29678.334

$ws.Range("H111").Value = 23644
$ws.Range("J111").Value = 23644
$ws.Range("L111").Value = 23644
$ws.Range("N111").Value = -31824

$ws.Range("H136").Value = 5027.1177
$ws.Range("I136").Value = 5162.2666
$ws.Range("J136").Value = 4013.5
$ws.Range("K136").Value = 15486.7998
$ws.Range("L136").Value = 12040.5
$ws.Range("M136").Value = -12936.7998
$ws.Range("N136").Value = -17140.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 365.66666
$ws.Range("J94").Value = 298.5
$ws.Range("L94").Value = 298.5
$ws.Range("N94").Value = -1200.5

$ws.Range("H107").Value = 2499.5
$ws.Range("I107").Value = 2499.5
$ws.Range("K107").Value = 2499.5
$ws.Range("M107").Value = -579.5

$ws.Range("H134").Value = 5831
$ws.Range("I134").Value = 5831
$ws.Range("K134").Value = 17493
$ws.Range("M134").Value = -14958

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1434.579
$ws.Range("I22").Value = 1020.44446
$ws.Range("J22").Value = 1807.3
$ws.Range("K22").Value = 1020.44446
$ws.Range("L22").Value = 1807.3
$ws.Range("M22").Value = -670.44446
$ws.Range("N22").Value = -2507.3

$ws.Range("H43").Value = 34128.5
$ws.Range("J43").Value = 34128.5
$ws.Range("L43").Value = 34128.5
$ws.Range("N43").Value = -34496.5

$ws.Range("H58").Value = 3717.5
$ws.Range("I58").Value = 3702
$ws.Range("K58").Value = 3702
$ws.Range("M58").Value = -3499

$ws.Range("H101").Value = 34128.5
$ws.Range("J101").Value = 34128.5
$ws.Range("L101").Value = 34128.5
$ws.Range("N101").Value = -40618.5

$ws.Range("H122").Value = 1448.6
$ws.Range("I122").Value = 1432.6666
$ws.Range("J122").Value = 1455.4286
$ws.Range("K122").Value = 4297.9998
$ws.Range("L122").Value = 4366.2858
$ws.Range("M122").Value = -1847.9998
$ws.Range("N122").Value = -9266.2858

$ws.Range("H136").Value = 3717.5
$ws.Range("I136").Value = 3702
$ws.Range("K136").Value = 11106
$ws.Range("M136").Value = -8556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5665
$ws.Range("I62").Value = 5665
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 16995
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -16309
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 5665
$ws.Range("I65").Value = 5665
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 50985
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -47553
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2202.8
$ws.Range("I80").Value = 1536
$ws.Range("K80").Value = 1536
$ws.Range("M80").Value = -538

$ws.Range("H83").Value = 2202.8
$ws.Range("I83").Value = 1536
$ws.Range("K83").Value = 7680
$ws.Range("M83").Value = -2688

$ws.Range("H102").Value = 3855.3635
$ws.Range("I102").Value = 4040.9
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 4040.9
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -2418.9
$ws.Range("N102").Value = -5244

$ws.Range("H104").Value = 49911.25
$ws.Range("J104").Value = 49911.25
$ws.Range("L104").Value = 49911.25
$ws.Range("N104").Value = -56899.25

$ws.Range("H122").Value = 4469.5713
$ws.Range("I122").Value = 3477.4
$ws.Range("K122").Value = 10432.2
$ws.Range("M122").Value = -7982.200000000001

$ws.Range("H132").Value = 3534.818
$ws.Range("I132").Value = 3534.818
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10604.454
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8074.454000000002
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2430.4
$ws.Range("I7").Value = 2283.1667
$ws.Range("J7").Value = 2651.25
$ws.Range("K7").Value = 2283.1667
$ws.Range("L7").Value = 2651.25
$ws.Range("M7").Value = -2171.1667
$ws.Range("N7").Value = -2875.25

$ws.Range("H95").Value = 65000
$ws.Range("J95").Value = 65000
$ws.Range("L95").Value = 65000
$ws.Range("N95").Value = -70492

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H126").Value = 2430.4
$ws.Range("I126").Value = 2283.1667
$ws.Range("J126").Value = 2651.25
$ws.Range("K126").Value = 6849.500100000001
$ws.Range("L126").Value = 7953.75
$ws.Range("M126").Value = -4379.500100000001
$ws.Range("N126").Value = -12893.75

$ws.Range("H132").Value = 2588.4443
$ws.Range("I132").Value = 2316.3333
$ws.Range("K132").Value = 6948.999899999999
$ws.Range("M132").Value = -4418.999899999999

$ws.Range("H136").Value = 55558508
$ws.Range("I136").Value = 3512.5715
$ws.Range("K136").Value = 10537.7145
$ws.Range("M136").Value = -7987.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 142858600
$ws.Range("I100").Value = 166667700
$ws.Range("K100").Value = 333335400
$ws.Range("M100").Value = -333334859

$ws.Range("H122").Value = 5313.143
$ws.Range("I122").Value = 5313.143
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15939.429
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -13489.429
$ws.Range("N122").ClearContents()

$ws.Range("H124").Value = 95000
$ws.Range("J124").Value = 95000
$ws.Range("L124").Value = 95000
$ws.Range("N124").Value = -104820

$ws.Range("H132").Value = 2075.6
$ws.Range("I132").Value = 2075.6
$ws.Range("K132").Value = 6226.799999999999
$ws.Range("M132").Value = -3696.799999999999

$ws.Range("H135").Value = 177143
$ws.Range("J135").Value = 177143
$ws.Range("L135").Value = 177143
$ws.Range("N135").Value = -187283

$ws.Range("H140").Value = 24659
$ws.Range("J140").Value = 24659
$ws.Range("L140").Value = 24659
$ws.Range("N140").Value = -35019

$ws.Range("H141").Value = 29999
$ws.Range("J141").Value = 29999
$ws.Range("L141").Value = 29999
$ws.Range("N141").Value = -40359

